# CheetahProcessing.xlsx - "Changes of config file"
#
# The commit replaces a batch of PackageTrackNum / ShipmentTrackNum values
# (columns C and D) on rows 2-22 with a fresh set of tracking numbers, and
# flips the Q3 pass/fail flag from "FAIL" to "Pass".
#
# Because these tracking numbers are strings of digits, a plain
# `.Value = "320017958830"` assignment would be auto-coerced to a number by
# Excel's type inference (losing the leading-text semantics / t="s" shared
# string storage the workbook originally used). To keep them as genuine
# text values we briefly force the cell to Text number format before the
# assignment, then restore the cell to the default "Normal" style so no
# lasting formatting change is left behind (matches the original cells,
# which carried no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

# Column C / D tracking-number replacements (row order matches the sheet).
Set-TextValue "C2"  "320017958830"
Set-TextValue "C3"  "320017958841"

# Q3: pass/fail indicator flips from FAIL to Pass (plain text, no coercion risk).
$ws.Range("Q3").Value = "Pass"

Set-TextValue "C4"  "320017958874"

Set-TextValue "C5"  "320017958896"
Set-TextValue "D5"  "320017958896"

Set-TextValue "C6"  "320017958933"
Set-TextValue "D6"  "320017958933"

Set-TextValue "C7"  "320017958955"
Set-TextValue "D7"  "320017958955"

Set-TextValue "C8"  "320017958988"
Set-TextValue "C9"  "320017959002"
Set-TextValue "C10" "320017959035"
Set-TextValue "C11" "320017959057"
Set-TextValue "C12" "320017959090"

Set-TextValue "C13" "320017959116"
Set-TextValue "D13" "320017959116"

Set-TextValue "C14" "320017959149"
Set-TextValue "D14" "320017959149"

Set-TextValue "C15" "320017959160"
Set-TextValue "D15" "320017959160"

Set-TextValue "C16" "320017959208"
Set-TextValue "D16" "320017959208"

Set-TextValue "C17" "320017959220"
Set-TextValue "D17" "320017959220"

Set-TextValue "C18" "320017959263"
Set-TextValue "C19" "320017959285"
Set-TextValue "C20" "320017959311"
Set-TextValue "C21" "320017959333"
Set-TextValue "C22" "320017959366"
